# Atualização automática de PALMARES_DO_SUL.xlsx
#
# - Renomeia a planilha "Paineis DARQ" para "PAINEIS DARQ"
# - Renomeia a planilha "Recolhimento x Eliminacao" para "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove a planilha "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

# Remove the "Desarquivamentos Pendentes" worksheet entirely.
$sheetToRemove = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$sheetToRemove.Delete() | Out-Null

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$painelSheet = $wb.Worksheets.Item("Paineis DARQ")
$painelSheet.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$recolhimentoSheet = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolhimentoSheet.Name = "RECOLHIMENTO X ELIMINAÇÃO"
